# Add DigiKey purchase links column + fix "Resister" -> "Resistor" typo on the BOM sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: fix spelling "Resister" -> "Resistor" for the two resistor rows ---
$ws.Range("A6").Value = "Resistor 10K 1206"
$ws.Range("A7").Value = "Resistor 470R 1206"

# --- Column D: new "DIGIKEY LINK" header + per-part purchase links ---
$ws.Range("D1").Value = "DIGIKEY LINK"
$ws.Range("D4").Value = "https://www.digikey.jp/en/products/detail/kemet/C1206C104K5RAC7210/3317012"
$ws.Range("D5").Value = "https://www.digikey.jp/en/products/detail/c&k/JS202011CQN/1640097?utm_adgroup=&utm_source=google&utm_medium=cpc&utm_campaign=PMax%20Shopping_Product_New%20Customer%20Acquisition&utm_term=&productid=1640097&utm_content=&utm_id=go_cmp-19897039674_adg-_ad-__dev-c_ext-_prd-1640097_sig-EAIaIQobChMIqub9taCMhAMVviN7Bx16jgQLEAQYAyABEgKK6_D_BwE&gad_source=1&gclid=EAIaIQobChMIqub9taCMhAMVviN7Bx16jgQLEAQYAyABEgKK6_D_BwE"
$ws.Range("D6").Value = "https://www.digikey.jp/en/products/detail/yageo/RC1206FR-0710KL/728483"
$ws.Range("D7").Value = "https://www.digikey.jp/en/products/detail/stackpole-electronics-inc/RMCF1206JT470R/1753845"
$ws.Range("D8").Value = "https://www.digikey.jp/en/products/detail/würth-elektronik/151033RS03000/4490003"
$ws.Range("D9").Value = "https://www.digikey.jp/en/products/detail/onsemi/NC7S04M5X-L22090/16910805"

# --- Selection moves from D4 to B7 ---
$ws.Range("B7").Select()
